$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row values (A: Subject, B: Helps, C: Problems, D: Day, E: Gain)
$rows = @{
    2  = @("Physics", "Understand forces—catapult design, castle defenses—enhancing military precision.", "Siege warfare inefficiencies; weak fortifications against dragons or trebuchets.", 10, "Stronger walls, fewer siege losses")
    3  = @("Chemistry", "Craft wildfire, fertilizers—control resources, boost agriculture, wield power.", "Famine from poor harvests; lack of alchemical defenses against foes like Cersei’s wildfire plots.", 15, "More crops, wildfire to deter invaders")
    4  = @("Mathematics", "Calculate taxes, troop logistics—ensure fair rule, efficient supply lines.", "Misallocated resources; inaccurate levies causing unrest among smallfolk.", 5, "Fair taxes, well-fed armies")
    5  = @("Discrete Mathematics", "Model alliances, voting systems—optimize diplomacy, council decisions.", "Feudal disputes; unclear succession lines fueling wars (e.g., War of the Five Kings).", 20, "Peaceful alliances, clear succession")
    6  = @("Engineering Mathematics", "Design bridges, aqueducts—improve infrastructure, unify the realm.", "Broken trade routes; flooded Riverlands disrupting food supply.", 25, "Better roads, steady food from aqueducts")
    7  = @("Digital Design", "Build signaling systems (e.g., raven relays)—speed communication across kingdoms.", "Slow message delivery; miscommunication sparking rebellions (e.g., Robb’s campaigns).", 30, "Faster news, fewer rebellions")
    8  = @("Computer Organisation and Architecture", "Structure data flow—organize royal records, troop movements digitally.", "Lost scrolls; chaotic command during battles like Blackwater.", 35, "Orderly records, swift battle commands")
    9  = @("Programming and Data Structures", "Code tools like ``raven_tally.py``—track resources, messages, automate rule.", "Inefficient raven messaging; untracked supplies during winter sieges.", 3, "Tracked ravens, stocked granaries")
    11 = @("Theory of Computation", "Predict system limits—ensure scalable governance as kingdoms grow.", "Overstretched rule post-war; failing to manage expanded territories after Daenerys’s conquests.", 50, "Stable rule over vast lands")
    12 = @("Compiler Design", "Translate decrees to code—standardize commands for maesters, lords.", "Misinterpreted royal orders; inconsistent law enforcement across regions.", 55, "Clear laws, uniform justice")
    13 = @("Operating Systems", "Manage kingdom processes—allocate resources, prioritize tasks like a king’s OS.", "Overlapping duties among lords; resource hoarding by Houses like Lannister.", 60, "Fair resource split, efficient tasks")
    14 = @("Database Management Systems", "Store folk’s data—track taxes, fealties, harvests in a royal database.", "Lost lineage records (e.g., Jon Snow’s claim); untracked grain stores leading to starvation.", 45, "Known lineage, full granaries")
    17 = @("Robotics", "Automate labor—forge golems for fields, walls—ease burdens, bolster defenses.", "Smallfolk exhaustion; crumbling defenses against invaders (e.g., Wildlings breaching the Wall).", 80, "Rested workers, unbreached walls")
    18 = @("Mechanical Engineering", "Build siege engines, water mills—strengthen war and peace efforts.", "Weak trebuchets failing at sieges; insufficient grain milling during winters.", 75, "Stronger sieges, milled grain aplenty")
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
    $ws.Cells.Item($r, 5).Value = $vals[4]
}
